$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B20").Value = 6169
$ws.Range("C20").Value = 979
$ws.Range("D20").Value = 5575743
$ws.Range("E20").Value = 903.8325498460042
$ws.Range("F20").Value = 6.564173432371745
$ws.Range("G20").Value = 3.927813163481964
$ws.Range("H20").Value = 26.12882699082244
